$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Averaged intensities now include the Gaussian-Quadrature scheme moved up next to
# "Ring Perpendicular to *" and three new spiral sampling schemes, with the remaining
# rotation/grid schemes shifted down to make room. Rows 17-19 are brand new.
$rows = @(
  @{ Row = 10; A = 8; B = "Gaussian-Quadrature"; Vals = @(1.02888972242382, 0.9467527183466591, 1.006267620948166, 0.9859838365847405, 1.02888972242382, 0.9467527183466591, 1.010908137415083, 0.9866228137006403, 1.006267620948166, 0.9641314439970358, 1.02888972242382, 0.9765101696474126, 0.9919734745758463, 0.9919779892955387) },
  @{ Row = 11; A = 9; B = "Spiral-90deg-10rot-5space"; Vals = @(1.063144878672589, 0.8718982202546992, 1.034318162922135, 0.9701860914520705, 1.063144878672589, 0.8718982202546992, 1.040294565969237, 0.9746022497151352, 1.021402246213574, 0.9157216013569897, 1.063144878672589, 0.9531081915884173, 0.9848868383253734, 0.9864460020695538) },
  @{ Row = 12; A = 10; B = "Spiral-90deg-15rot-5space"; Vals = @(1.062735922220136, 0.8723433355314836, 1.034269962343455, 0.970293601554107, 1.062735922220136, 0.8723433355314836, 1.0401372217318, 0.9747237756893432, 1.021247363968429, 0.9160226533190126, 1.062735922220136, 0.9533066489374693, 0.9849107054122952, 0.9864717295447207) },
  @{ Row = 13; A = 11; B = "Spiral-90deg-10rot-3space"; Vals = @(1.06306815280823, 0.8719509310398167, 1.034354012378342, 0.970208129526667, 1.06306815280823, 0.8719509310398167, 1.040263894517532, 0.9746304761116122, 1.021344916962991, 0.9157765967129831, 1.06306815280823, 0.9531524717090791, 0.984895306438264, 0.9864496387572718) },
  @{ Row = 14; A = 12; B = "NoRotation-tilt60deg"; Vals = @(1.064012000000002, 0.8549680000000002, 1.037903999999999, 0.9643639999999996, 1.064012000000002, 0.8549680000000002, 1.048971999999999, 0.9747720000000005, 1.028563999999999, 0.9062240000000009, 1.064012000000002, 0.9464359999999995, 0.9803120000000001, 0.9849724999999999) },
  @{ Row = 15; A = 13; B = "Rotation-NoTilt"; Vals = @(1.11, 0.7944500000000024, 1.043050000000001, 0.9475124999999988, 1.11, 0.7944500000000024, 1.07, 0.96, 1.05, 0.87, 1.11, 0.9187500000000017, 0.9737531250000006, 0.9806265625000004) },
  @{ Row = 16; A = 14; B = "Rotation-60detTilt"; Vals = @(1.062597072486396, 0.8777196858368004, 1.023768468889597, 0.9678807928832042, 1.062597072486396, 0.8777196858368004, 1.038550623744004, 0.9743610972160013, 1.026906695065596, 0.9224142880768027, 1.062597072486396, 0.9507440773631989, 0.9829915050239995, 0.9867748405248002) },
  @{ Row = 17; A = 15; B = "HexGrid-90degTilt5degRes"; Vals = @(0.9959084322891807, 0.9950620489413309, 0.9949473699897788, 0.9951677999633605, 0.9959084322891807, 0.9950620489413309, 0.9950734519140056, 0.9955598358402616, 0.9954845923855136, 0.9950443528309405, 0.9959084322891807, 0.9950047094655549, 0.9952714127959127, 0.9952809855192966) },
  @{ Row = 18; A = 16; B = "HexGrid-90degTilt22p5degRes"; Vals = @(0.989452046119333, 1.001925943980766, 0.9955028962436309, 0.9963759342021271, 0.989452046119333, 1.001925943980766, 0.9927013703325483, 0.9973050308813205, 0.9924674077872336, 0.9995685899046446, 0.989452046119333, 0.9987144201121982, 0.9958142051364641, 0.9956624024314504) },
  @{ Row = 19; A = 17; B = "HexGrid-60degTilt5degRes"; Vals = @(0.9800078939874022, 1.016415740110491, 0.9905726222273578, 1.00015042395527, 0.9800078939874022, 1.016415740110491, 0.986881509433071, 0.9994405788272167, 0.9900236527316723, 1.009608526493999, 0.9800078939874022, 1.003494181168924, 0.9967866700701302, 0.9966376184708099) }
)

foreach ($r in $rows) {
  $ws.Cells.Item($r.Row, 1).Value = $r.A
  $ws.Cells.Item($r.Row, 2).Value = $r.B
  $col = 3
  foreach ($v in $r.Vals) {
    $ws.Cells.Item($r.Row, $col).Value = $v
    $col++
  }
}

# New rows 17-19 need the same "index column" style (bold, centered, bordered)
# already used by column A in rows 2-16.
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17:A19").PasteSpecial(-4122) | Out-Null
